$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1939.1111
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 1931.5
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 1931.5
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -2583.5
$ws.Range("H116").Value = 4291.4736
$ws.Range("I116").Value = 5317.5557
$ws.Range("J116").Value = 3368
$ws.Range("K116").Value = 5317.5557
$ws.Range("L116").Value = 3368
$ws.Range("M116").Value = -1875.5557
$ws.Range("N116").Value = -10252
$ws.Range("H135").Value = 1513.7931
$ws.Range("I135").Value = 1550
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 13950
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -11415
$ws.Range("N135").Value = -15870
$ws.Range("H137").Value = 3079.257
$ws.Range("I137").Value = 1431.091
$ws.Range("J137").Value = 5868.4614
$ws.Range("K137").Value = 4293.272999999999
$ws.Range("L137").Value = 17605.3842
$ws.Range("M137").Value = -1743.272999999999
$ws.Range("N137").Value = -22705.3842
$ws.Range("H138").Value = 5557986
$ws.Range("I138").Value = 1214.5938
$ws.Range("J138").Value = 19236192
$ws.Range("K138").Value = 3643.7814
$ws.Range("L138").Value = 57708576
$ws.Range("M138").Value = 1496.2186
$ws.Range("N138").Value = -57718856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 997.2069
$ws.Range("I2").Value = 906.1429000000001
$ws.Range("J2").Value = 1236.25
$ws.Range("K2").Value = 906.1429000000001
$ws.Range("L2").Value = 1236.25
$ws.Range("M2").Value = -793.1429000000001
$ws.Range("N2").Value = -1462.25
$ws.Range("H32").Value = 7868.9
$ws.Range("I32").Value = 9435.462
$ws.Range("K32").Value = 9435.462
$ws.Range("M32").Value = -9148.462
$ws.Range("H61").Value = 3224.8635
$ws.Range("I61").Value = 3052.611
$ws.Range("K61").Value = 3052.611
$ws.Range("M61").Value = -2840.611
$ws.Range("H74").Value = 2158.9707
$ws.Range("I74").Value = 1356.3214
$ws.Range("J74").Value = 5904.6665
$ws.Range("K74").Value = 1356.3214
$ws.Range("L74").Value = 5904.6665
$ws.Range("M74").Value = -482.3214
$ws.Range("N74").Value = -7652.6665
$ws.Range("H77").Value = 2158.9707
$ws.Range("I77").Value = 1356.3214
$ws.Range("J77").Value = 5904.6665
$ws.Range("K77").Value = 6781.607
$ws.Range("L77").Value = 29523.3325
$ws.Range("M77").Value = -2413.607
$ws.Range("N77").Value = -38259.3325
$ws.Range("H116").Value = 997.2069
$ws.Range("I116").Value = 906.1429000000001
$ws.Range("J116").Value = 1236.25
$ws.Range("K116").Value = 906.1429000000001
$ws.Range("L116").Value = 1236.25
$ws.Range("M116").Value = 1387.8571
$ws.Range("N116").Value = -5824.25
$ws.Range("H136").Value = 3224.8635
$ws.Range("I136").Value = 3052.611
$ws.Range("K136").Value = 9157.832999999999
$ws.Range("M136").Value = -6607.832999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 997.2069
$ws.Range("I3").Value = 906.1429000000001
$ws.Range("J3").Value = 1236.25
$ws.Range("K3").Value = 906.1429000000001
$ws.Range("L3").Value = 1236.25
$ws.Range("M3").Value = -792.1429000000001
$ws.Range("N3").Value = -1464.25
$ws.Range("H99").Value = 1411.3334
$ws.Range("J99").Value = 2196
$ws.Range("L99").Value = 2196
$ws.Range("N99").Value = -5192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 137.06897
$ws.Range("I7").Value = 80.052635
$ws.Range("J7").Value = 245.4
$ws.Range("K7").Value = 80.052635
$ws.Range("L7").Value = 245.4
$ws.Range("M7").Value = 32.947365
$ws.Range("N7").Value = -471.4
$ws.Range("H31").Value = 777484.1
$ws.Range("I31").Value = 10607.071
$ws.Range("J31").Value = 1603351.8
$ws.Range("K31").Value = 10607.071
$ws.Range("L31").Value = 1603351.8
$ws.Range("M31").Value = -10312.071
$ws.Range("N31").Value = -1603941.8
$ws.Range("H34").Value = 777484.1
$ws.Range("I34").Value = 10607.071
$ws.Range("J34").Value = 1603351.8
$ws.Range("K34").Value = 10607.071
$ws.Range("L34").Value = 1603351.8
$ws.Range("M34").Value = -10405.071
$ws.Range("N34").Value = -1603755.8
$ws.Range("H99").Value = 1223.0714
$ws.Range("I99").Value = 1192.3
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1192.3
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = 305.7
$ws.Range("N99").Value = -4296
$ws.Range("H121").Value = 31740.8
$ws.Range("I121").Value = 15000
$ws.Range("J121").Value = 33600.89
$ws.Range("K121").Value = 15000
$ws.Range("L121").Value = 33600.89
$ws.Range("M121").Value = -13690
$ws.Range("N121").Value = -36220.89
$ws.Range("H126").Value = 1223.0714
$ws.Range("I126").Value = 1192.3
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 3576.9
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1106.9
$ws.Range("N126").Value = -8840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 205
$ws.Range("I23").Value = 230
$ws.Range("J23").Value = 202.22223
$ws.Range("K23").Value = 690
$ws.Range("L23").Value = 606.66669
$ws.Range("M23").Value = -455
$ws.Range("N23").Value = -1076.66669
$ws.Range("H122").Value = 2238.08
$ws.Range("I122").Value = 2089.8333
$ws.Range("K122").Value = 18808.4997
$ws.Range("M122").Value = -16358.4997
$ws.Range("H132").Value = 1901.6364
$ws.Range("I132").Value = 902.6667
$ws.Range("J132").Value = 2276.25
$ws.Range("K132").Value = 8124.0003
$ws.Range("L132").Value = 20486.25
$ws.Range("M132").Value = -5594.0003
$ws.Range("N132").Value = -25546.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2498.8096
$ws.Range("I16").Value = 1808.4375
$ws.Range("J16").Value = 4708
$ws.Range("K16").Value = 1808.4375
$ws.Range("L16").Value = 4708
$ws.Range("M16").Value = -1638.4375
$ws.Range("N16").Value = -5048
$ws.Range("H132").Value = 14294347
$ws.Range("I132").Value = 5077.067
$ws.Range("J132").Value = 25011300
$ws.Range("K132").Value = 15231.201
$ws.Range("L132").Value = 75033900
$ws.Range("M132").Value = -12701.201
$ws.Range("N132").Value = -75038960
$ws.Range("H136").Value = 34489916
$ws.Range("I136").Value = 55558010
$ws.Range("J136").Value = 14856.363
$ws.Range("K136").Value = 166674030
$ws.Range("L136").Value = 44569.089
$ws.Range("M136").Value = -166671480
$ws.Range("N136").Value = -49669.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2130.3333
$ws.Range("I107").Value = 2898
$ws.Range("J107").Value = 595
$ws.Range("K107").Value = 8694
$ws.Range("L107").Value = 1785
$ws.Range("M107").Value = -6774
$ws.Range("N107").Value = -5625
$ws.Range("H136").Value = 3847800.2
$ws.Range("I136").Value = 5557045.5
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 16671136.5
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -16668586.5
$ws.Range("N136").Value = -11094
